# Remove the "LKN " prefix (as in "LKN AA.00.0020") that was accidentally left
# inside some of the question texts in column B ("Frage").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $text = $cell.Value2
    if ($text -ne $null -and $text -like "*LKN *") {
        $cell.Value = $text.Replace("LKN ", "")
    }
}

$ws.Range("B72").Select()
